# Deploying to gh-pages from @ IDGRLP/Tumorkonferenzen-IG@a6ec8440d71eb4ead340661347be820f66a9ac52 🚀
#
# - Rename the "Include from Histopathologisc" sheet to "Include #0".
# - Metadata sheet: insert a new "Jurisdiction" property row (with an
#   empty value) right after "Contact", pushing Description/Purpose/
#   Copyright/Immutable down by one row.
# - Metadata sheet: refresh the "Date" property's value.

$wb = $excel.ActiveWorkbook

# --- Rename the second worksheet -------------------------------------------------
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

# --- Metadata worksheet -----------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Extend the body formatting (same style as the existing property rows, e.g. row 10)
# down through the new last row (15) *before* we write any values there, so the
# freshly-touched cells pick up the existing "s=2" look instead of Excel's
# bare default style.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B15").PasteSpecial(-4122)

# Shift the existing "Description"..."Immutable" rows down by one (row 11 -> 12,
# 12 -> 13, 13 -> 14, 14 -> 15) to make room for the new "Jurisdiction" row.
for ($i = 14; $i -ge 11; $i--) {
    $dst = $i + 1
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($i, 1).Value2
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($i, 2).Value2
}

# New row 11: Jurisdiction property with an empty value.
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Refresh the Date property's value (row 8, column B).
$ws.Cells.Item(8, 2).Value = "2024-09-17T19:55:11+00:00"
